$wb = $excel.ActiveWorkbook

# --- Sheet2: bump the schedule-number counter and data-recorder index ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B1").Value = 118
$ws2.Range("B2").Value = 24

# --- Sheet3: append a new tracking row (row 23) ---
$ws3 = $wb.Worksheets.Item("Sheet3")

# Columns A, B and D hold digit-only strings that must stay text (shared
# strings), not numbers, so force text entry via the leading apostrophe and
# then drop the resulting "quote prefix" style off the cell itself.
$ws3.Range("A23").Value = "'3013696628"
$ws3.Range("A23").Style = "Normal"

$ws3.Range("B23").Value = "'1000004650"
$ws3.Range("B23").Style = "Normal"

# Column C reuses the existing "schedNum" label text already used above.
$ws3.Range("C23").Value = "schedNum"

$ws3.Range("D23").Value = "'13188908"
$ws3.Range("D23").Style = "Normal"
